$wb = $excel.ActiveWorkbook

# Add new worksheet right after the existing "IF" sheet
$ifSheet = $wb.Worksheets.Item("IF")
$newSheet = $wb.Worksheets.Add($null, $ifSheet)
$newSheet.Name = "Sheet1"

# Populate the new sheet's data
$newSheet.Range("C4").Value = "ExcelWriter"
$newSheet.Range("C5").Value = "　・WorkBook"

$newSheet.Range("E4").Value = "OneResultWriteProcedure"
$newSheet.Range("E5").Value = "　・Sheet"

$newSheet.Range("H4").Value = "OperationCellUtil"
$newSheet.Range("H5").Value = "　・OpeCell"
$newSheet.Range("H6").Value = "　・initXpositon"
$newSheet.Range("H7").Value = "　・initYpositon"

$newSheet.Range("E9").Value = "OperatableCell"
$newSheet.Range("E10").Value = "　"

$newSheet.Range("E18").Value = "OpeCell"
$newSheet.Range("E19").Value = "　・x"
$newSheet.Range("E20").Value = "　・y"
$newSheet.Range("E21").Value = "　・getX"
$newSheet.Range("E22").Value = "　・getY"
$newSheet.Range("E23").Value = "　・incrementX"
$newSheet.Range("E24").Value = "　・incrementY"

# Select E10 as the active cell on the new sheet
$newSheet.Range("E10").Select()

# Make the new sheet the active tab
$newSheet.Activate()
